# Add export for large publications + started to add weight_type (for Scores).
#
# On the "Scores" worksheet, insert a new column "Type of Variant Weight"
# right after "Number of Interaction Terms" (new column K), pushing the
# existing "PGS Publication (PGP) ID" ... "License/Terms of Use" columns
# one slot to the right (old K:S -> new L:T). Populate the header + the
# single data row, and re-point the FTP-link hyperlink (which lived on the
# old column R, now column S) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scores")

# Remember the hyperlink target before we shuffle columns around.
$ftpLink = "http://ftp.ebi.ac.uk/pub/databases/spot/pgs/scores/PGS1/ScoringFiles/PGS1.txt.gz"

# Insert a new column at K - everything from old K onward shifts right by one.
$ws.Columns.Item(11).Insert()

# New column header + value.
$ws.Range("K1").Value = "Type of Variant Weight"
$ws.Range("K2").Value = "log(OR)"

# The hyperlink that used to anchor at R2 now needs to point at S2 (it did
# not move automatically with the column insert).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("S2"), $ftpLink)
$ws.Range("S2").Style = "Hyperlink"

Write-Output "Scores sheet updated with Type of Variant Weight column"
